$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.02648366666666667
$ws.Range("N2").Value = 0.07945099999999999
$ws.Range("O2").Value = 0.001430039273477916
$ws.Range("P2").Value = 0.001430039273477917
$ws.Range("Q2").Value = 0.022061556425
$ws.Range("R2").Value = 0.198554007825
$ws.Range("S2").Value = 0.001430039273477916
$ws.Range("T2").Value = 0.001430039273477917

# Row 3 updates
$ws.Range("O3").Value = 0.7016741634339546
$ws.Range("P3").Value = 0.7016741634339547
$ws.Range("S3").Value = 0.7016741634339546
$ws.Range("T3").Value = 0.7016741634339547

# Row 4 updates
$ws.Range("O4").Value = 0.2968957972925674
$ws.Range("P4").Value = 0.2968957972925675
$ws.Range("S4").Value = 0.2968957972925674
$ws.Range("T4").Value = 0.2968957972925675
